$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-05-03 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-05-04 Saturday", 2)

# Update the multiplication problems in the table. Cells are addressed by
# (row, column) rather than by text search because some of the new values
# collide with old values used elsewhere in the table (e.g. "71×20=" is both
# a new value and an old value), which would make a plain text find/replace
# unsafe and order-dependent.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "75×83="
$t.Cell(1, 2).Range.Text  = "82×16="
$t.Cell(1, 3).Range.Text  = "53×85="
$t.Cell(1, 4).Range.Text  = "75×56="
$t.Cell(1, 5).Range.Text  = "55×46="

$t.Cell(5, 1).Range.Text  = "38×64="
$t.Cell(5, 2).Range.Text  = "71×20="
$t.Cell(5, 3).Range.Text  = "36×46="
$t.Cell(5, 4).Range.Text  = "77×94="
$t.Cell(5, 5).Range.Text  = "62×27="

$t.Cell(10, 1).Range.Text = "28×61="
$t.Cell(10, 2).Range.Text = "46×51="
$t.Cell(10, 3).Range.Text = "73×75="
$t.Cell(10, 4).Range.Text = "94×24="
$t.Cell(10, 5).Range.Text = "23×87="

$t.Cell(15, 1).Range.Text = "84×22="
$t.Cell(15, 2).Range.Text = "43×25="
$t.Cell(15, 3).Range.Text = "14×81="
$t.Cell(15, 4).Range.Text = "87×89="
$t.Cell(15, 5).Range.Text = "67×75="

$t.Cell(20, 1).Range.Text = "24×31="
$t.Cell(20, 2).Range.Text = "63×67="
$t.Cell(20, 3).Range.Text = "93×27="
$t.Cell(20, 4).Range.Text = "27×48="
$t.Cell(20, 5).Range.Text = "45×44="
